$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Clear stray summary-stat columns from rows 3-5 (unique/top/freq) that
# belonged to the wrong rows.
$ws.Cells.Item(3, 7).ClearContents()
$ws.Cells.Item(3, 8).ClearContents()
$ws.Cells.Item(3, 9).ClearContents()
$ws.Cells.Item(3, 10).ClearContents()
$ws.Cells.Item(3, 11).ClearContents()
$ws.Cells.Item(3, 12).ClearContents()
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).ClearContents()
$ws.Cells.Item(3, 15).ClearContents()
$ws.Cells.Item(3, 16).ClearContents()
$ws.Cells.Item(3, 17).ClearContents()
$ws.Cells.Item(3, 18).ClearContents()
$ws.Cells.Item(3, 19).ClearContents()
$ws.Cells.Item(3, 20).ClearContents()
$ws.Cells.Item(3, 21).ClearContents()
$ws.Cells.Item(3, 22).ClearContents()
$ws.Cells.Item(3, 23).ClearContents()
$ws.Cells.Item(3, 25).ClearContents()
$ws.Cells.Item(3, 29).ClearContents()
$ws.Cells.Item(3, 33).ClearContents()

$ws.Cells.Item(4, 7).ClearContents()
$ws.Cells.Item(4, 8).ClearContents()
$ws.Cells.Item(4, 9).ClearContents()
$ws.Cells.Item(4, 10).ClearContents()
$ws.Cells.Item(4, 11).ClearContents()
$ws.Cells.Item(4, 12).ClearContents()
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(4, 15).ClearContents()
$ws.Cells.Item(4, 16).ClearContents()
$ws.Cells.Item(4, 17).ClearContents()
$ws.Cells.Item(4, 18).ClearContents()
$ws.Cells.Item(4, 19).ClearContents()
$ws.Cells.Item(4, 20).ClearContents()
$ws.Cells.Item(4, 21).ClearContents()
$ws.Cells.Item(4, 22).ClearContents()
$ws.Cells.Item(4, 23).ClearContents()
$ws.Cells.Item(4, 25).ClearContents()
$ws.Cells.Item(4, 29).ClearContents()
$ws.Cells.Item(4, 33).ClearContents()

$ws.Cells.Item(5, 7).ClearContents()
$ws.Cells.Item(5, 8).ClearContents()
$ws.Cells.Item(5, 9).ClearContents()
$ws.Cells.Item(5, 10).ClearContents()
$ws.Cells.Item(5, 11).ClearContents()
$ws.Cells.Item(5, 12).ClearContents()
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(5, 15).ClearContents()
$ws.Cells.Item(5, 16).ClearContents()
$ws.Cells.Item(5, 17).ClearContents()
$ws.Cells.Item(5, 18).ClearContents()
$ws.Cells.Item(5, 19).ClearContents()
$ws.Cells.Item(5, 20).ClearContents()
$ws.Cells.Item(5, 21).ClearContents()
$ws.Cells.Item(5, 22).ClearContents()
$ws.Cells.Item(5, 23).ClearContents()
$ws.Cells.Item(5, 25).ClearContents()
$ws.Cells.Item(5, 29).ClearContents()
$ws.Cells.Item(5, 33).ClearContents()

# Populate rows 6-12 (mean/std/min/25%/50%/75%/max) with the correct
# per-column descriptive statistics.
# Row 6
$ws.Cells.Item(6, 7).Value = 40.60432299793592
$ws.Cells.Item(6, 8).Value = 26.44206943848964
$ws.Cells.Item(6, 9).Value = 42.81318016591187
$ws.Cells.Item(6, 10).Value = 49.4768054228887
$ws.Cells.Item(6, 11).Value = 1566.752522467704
$ws.Cells.Item(6, 12).Value = 1653.625127759921
$ws.Cells.Item(6, 13).Value = 2425.883965411414
$ws.Cells.Item(6, 14).Value = 124379.4953272435
$ws.Cells.Item(6, 15).Value = 9.111649313701808
$ws.Cells.Item(6, 16).Value = 37.97458256181991
$ws.Cells.Item(6, 17).Value = 64.48632034329007
$ws.Cells.Item(6, 18).Value = 1.224350214138081
$ws.Cells.Item(6, 19).Value = 1.08688534414102
$ws.Cells.Item(6, 20).Value = -2.897813695950064
$ws.Cells.Item(6, 21).Value = 49.69170728029123
$ws.Cells.Item(6, 22).Value = 3925786798.274295
$ws.Cells.Item(6, 23).Value = 45.08999531772432
$ws.Cells.Item(6, 25).Value = 0.08688650894950616
$ws.Cells.Item(6, 29).Value = 820.1013819865034
$ws.Cells.Item(6, 33).Value = 29.22209177644399

# Row 7
$ws.Cells.Item(7, 7).Value = 21.02235277283804
$ws.Cells.Item(7, 8).Value = 26.41242032004845
$ws.Cells.Item(7, 9).Value = 22.62590860512981
$ws.Cells.Item(7, 10).Value = 22.78286139374379
$ws.Cells.Item(7, 11).Value = 2032.941644460481
$ws.Cells.Item(7, 12).Value = 2021.725719384686
$ws.Cells.Item(7, 13).Value = 2017.510936805265
$ws.Cells.Item(7, 14).Value = 160183.2908472984
$ws.Cells.Item(7, 15).Value = 23.05367796400245
$ws.Cells.Item(7, 16).Value = 351.7387432869111
$ws.Cells.Item(7, 17).Value = 521.1813266575277
$ws.Cells.Item(7, 18).Value = 0.699713089277471
$ws.Cells.Item(7, 19).Value = 0.5992760263355361
$ws.Cells.Item(7, 20).Value = 32.53274329231496
$ws.Cells.Item(7, 21).Value = 26.38298018277067
$ws.Cells.Item(7, 22).Value = 13096910101.57792
$ws.Cells.Item(7, 23).Value = 25.87297132514584
$ws.Cells.Item(7, 25).Value = 0.3675053469566921
$ws.Cells.Item(7, 29).Value = 21714.27698207647
$ws.Cells.Item(7, 33).Value = 397.5019453200886

# Row 8
$ws.Cells.Item(8, 7).Value = 2.16036717202278
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 1.02448250499106
$ws.Cells.Item(8, 10).Value = 1.03274559193954
$ws.Cells.Item(8, 11).Value = 0
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 6.351908046707315
$ws.Cells.Item(8, 14).Value = 10.08283610702105
$ws.Cells.Item(8, 15).Value = 0.0015
$ws.Cells.Item(8, 16).Value = 0.008564348508246999
$ws.Cells.Item(8, 17).Value = 0.07263964552878199
$ws.Cells.Item(8, 18).Value = -1.4452784998966
$ws.Cells.Item(8, 19).Value = -0.258896562661341
$ws.Cells.Item(8, 20).Value = -842.151252707
$ws.Cells.Item(8, 21).Value = 0
$ws.Cells.Item(8, 22).Value = 1866576.53527322
$ws.Cells.Item(8, 23).Value = 0
$ws.Cells.Item(8, 25).Value = 0.00006063618290258458
$ws.Cells.Item(8, 29).Value = 0.0006432620052380001
$ws.Cells.Item(8, 33).Value = 0.0001714480194709703

# Row 9
$ws.Cells.Item(9, 7).Value = 23.48044685127255
$ws.Cells.Item(9, 8).Value = 2.17108773862599
$ws.Cells.Item(9, 9).Value = 25.17719770829013
$ws.Cells.Item(9, 10).Value = 31.34766011627065
$ws.Cells.Item(9, 11).Value = 35.61502562195389
$ws.Cells.Item(9, 12).Value = 54.84232474422455
$ws.Cells.Item(9, 13).Value = 786.8455343023638
$ws.Cells.Item(9, 14).Value = 12945.50886007589
$ws.Cells.Item(9, 15).Value = 0.7
$ws.Cells.Item(9, 16).Value = 6.844607536320606
$ws.Cells.Item(9, 17).Value = 10.71958658975
$ws.Cells.Item(9, 18).Value = 0.83862821426832
$ws.Cells.Item(9, 19).Value = 0.6661071171265258
$ws.Cells.Item(9, 20).Value = -5.6834794496
$ws.Cells.Item(9, 21).Value = 31.0969492849966
$ws.Cells.Item(9, 22).Value = 198630076.2029462
$ws.Cells.Item(9, 23).Value = 27.3476192806246
$ws.Cells.Item(9, 25).Value = 0.03027429419420288
$ws.Cells.Item(9, 29).Value = 1.159136290676698
$ws.Cells.Item(9, 33).Value = 1.375437490318808

# Row 10
$ws.Cells.Item(10, 7).Value = 37.94902687367665
$ws.Cells.Item(10, 8).Value = 18.56011789730495
$ws.Cells.Item(10, 9).Value = 40.70015947858415
$ws.Cells.Item(10, 10).Value = 48.56598142226279
$ws.Cells.Item(10, 11).Value = 622.5531833710047
$ws.Cells.Item(10, 12).Value = 746.8600732776019
$ws.Cells.Item(10, 13).Value = 1899.320061502804
$ws.Cells.Item(10, 14).Value = 54651.48072600651
$ws.Cells.Item(10, 15).Value = 2.7190094475
$ws.Cells.Item(10, 16).Value = 12.9509082136824
$ws.Cells.Item(10, 17).Value = 19.1603426405092
$ws.Cells.Item(10, 18).Value = 1.17109985413391
$ws.Cells.Item(10, 19).Value = 1.00882716143141
$ws.Cells.Item(10, 20).Value = 3.325495419
$ws.Cells.Item(10, 21).Value = 42.7114602373328
$ws.Cells.Item(10, 22).Value = 637823291.1727235
$ws.Cells.Item(10, 23).Value = 38.8014765635393
$ws.Cells.Item(10, 25).Value = 0.05219115266495034
$ws.Cells.Item(10, 29).Value = 3.011525833169725
$ws.Cells.Item(10, 33).Value = 3.387357080809172

# Row 11
$ws.Cells.Item(11, 7).Value = 55.8085593657264
$ws.Cells.Item(11, 8).Value = 44.71878761674657
$ws.Cells.Item(11, 9).Value = 59.69642665709583
$ws.Cells.Item(11, 10).Value = 67.430836296219
$ws.Cells.Item(11, 11).Value = 2443.749869229816
$ws.Cells.Item(11, 12).Value = 2709.88651757683
$ws.Cells.Item(11, 13).Value = 3471.141985412267
$ws.Cells.Item(11, 14).Value = 173821.0829215753
$ws.Cells.Item(11, 15).Value = 7.498422
$ws.Cells.Item(11, 16).Value = 23.7289934845
$ws.Cells.Item(11, 17).Value = 33.0313253835
$ws.Cells.Item(11, 18).Value = 1.59915043957711
$ws.Cells.Item(11, 19).Value = 1.44171485602766
$ws.Cells.Item(11, 20).Value = 8.438454714300001
$ws.Cells.Item(11, 21).Value = 61.7452197879924
$ws.Cells.Item(11, 22).Value = 2200145121.368382
$ws.Cells.Item(11, 23).Value = 55.8491649443871
$ws.Cells.Item(11, 25).Value = 0.09328730866713247
$ws.Cells.Item(11, 29).Value = 10.34662760911812
$ws.Cells.Item(11, 33).Value = 8.395804536415389

# Row 12
$ws.Cells.Item(12, 7).Value = 90.8780085125484
$ws.Cells.Item(12, 8).Value = 96.3711567124062
$ws.Cells.Item(12, 9).Value = 96.9258279083044
$ws.Cells.Item(12, 10).Value = 98.1338048317139
$ws.Cells.Item(12, 11).Value = 8935.80148567679
$ws.Cells.Item(12, 12).Value = 8780.787176347189
$ws.Cells.Item(12, 13).Value = 8892.806589289692
$ws.Cells.Item(12, 14).Value = 750544.426426754
$ws.Cells.Item(12, 15).Value = 290.72
$ws.Cells.Item(12, 16).Value = 13907.3775699265
$ws.Cells.Item(12, 17).Value = 16491.8032786885
$ws.Cells.Item(12, 18).Value = 4.88895043111182
$ws.Cells.Item(12, 19).Value = 3.77451006704652
$ws.Cells.Item(12, 20).Value = 177.6931560046
$ws.Cells.Item(12, 21).Value = 245.309152184057
$ws.Cells.Item(12, 22).Value = 157485039972.121
$ws.Cells.Item(12, 23).Value = 265.883848205148
$ws.Cells.Item(12, 25).Value = 13.76658700246231
$ws.Cells.Item(12, 29).Value = 1010044.90447419
$ws.Cells.Item(12, 33).Value = 15444.73098747262
